# Insert a new price record for "Feria Lagunitas de Puerto Montt - Cebollín"
# as row 279, pushing the existing rows 279:349 down to 280:350 and growing
# the used range from A1:R349 to A1:R350.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 279:349 down one row (this also carries the D-column's
# date number-format down into the newly created blank row 279).
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(279, 1).Value  = 4
$ws.Cells.Item(279, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(279, 3).Value  = "Los Lagos"
$ws.Cells.Item(279, 4).Value  = 44855
$ws.Cells.Item(279, 5).Value  = 10
$ws.Cells.Item(279, 6).Value  = 100112037
$ws.Cells.Item(279, 7).Value  = "Cebollín"
$ws.Cells.Item(279, 8).Value  = "Sin especificar"
$ws.Cells.Item(279, 9).Value  = "Primera"
$ws.Cells.Item(279, 10).Value = 180
$ws.Cells.Item(279, 11).Value = 6000
$ws.Cells.Item(279, 12).Value = 6000
$ws.Cells.Item(279, 13).Value = 6000
$ws.Cells.Item(279, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(279, 15).Value = "Región Metropolitana"
$ws.Cells.Item(279, 16).Value = 167
$ws.Cells.Item(279, 17).Value = 36
$ws.Cells.Item(279, 18).Value = "Hortaliza"
